# Auto-generated edit script applying cell-level corrections to the
# final_timetable workbook (mon/tue/wed/thur/fri sheets), per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("mon")
$ws.Cells.Item(2, 10).Value = ""  # J2: 'CSC442' -> None
$ws.Cells.Item(4, 4).Value = ""  # D4: 'CSC423' -> None
$ws.Cells.Item(4, 7).Value = "CSC424"  # G4: None -> 'CSC424'
$ws.Cells.Item(11, 3).Value = "CSC423"  # C11: None -> 'CSC423'
$ws.Cells.Item(13, 7).Value = "CSC425"  # G13: None -> 'CSC425'
$ws.Cells.Item(14, 7).Value = "CSC111"  # G14: None -> 'CSC111'
$ws.Cells.Item(15, 8).Value = ""  # H15: 'CSC425' -> None
$ws.Cells.Item(15, 9).Value = ""  # I15: 'CSC442' -> None
$ws.Cells.Item(20, 2).Value = "DLD221"  # B20: None -> 'DLD221'
$ws.Cells.Item(20, 3).Value = "DLD221"  # C20: None -> 'DLD221'
$ws.Cells.Item(20, 7).Value = "PHY111"  # G20: 'MAT111' -> 'PHY111'
$ws.Cells.Item(20, 8).Value = "PHY111"  # H20: None -> 'PHY111'
$ws.Cells.Item(21, 4).Value = ""  # D21: 'CST111' -> None
$ws.Cells.Item(21, 8).Value = ""  # H21: 'CHM111' -> None
$ws.Cells.Item(21, 9).Value = "CIT111"  # I21: 'CHM111' -> 'CIT111'
$ws.Cells.Item(21, 10).Value = "CIT111"  # J21: None -> 'CIT111'
$ws.Cells.Item(26, 8).Value = "CSC441"  # H26: None -> 'CSC441'
$ws.Cells.Item(26, 9).Value = "CSC441"  # I26: None -> 'CSC441'

$ws = $wb.Worksheets.Item("tue")
$ws.Cells.Item(4, 10).Value = ""  # J4: 'CSC423' -> None
$ws.Cells.Item(4, 11).Value = ""  # K4: 'CSC423' -> None
$ws.Cells.Item(14, 10).Value = "CSC424"  # J14: None -> 'CSC424'
$ws.Cells.Item(16, 5).Value = ""  # E16: 'CSC111' -> None
$ws.Cells.Item(17, 10).Value = "CSC423"  # J17: None -> 'CSC423'
$ws.Cells.Item(17, 11).Value = "CSC423"  # K17: None -> 'CSC423'
$ws.Cells.Item(20, 2).Value = ""  # B20: 'MAT111' -> None
$ws.Cells.Item(20, 3).Value = ""  # C20: 'MAT111' -> None
$ws.Cells.Item(20, 8).Value = "TMC421"  # H20: None -> 'TMC421'
$ws.Cells.Item(21, 4).Value = ""  # D21: 'CST111' -> None
$ws.Cells.Item(21, 5).Value = ""  # E21: 'GST111' -> None
$ws.Cells.Item(24, 9).Value = ""  # I24: 'CSC441' -> None
$ws.Cells.Item(24, 10).Value = ""  # J24: 'CSC441' -> None

$ws = $wb.Worksheets.Item("wed")
$ws.Cells.Item(3, 9).Value = ""  # I3: 'CSC424' -> None
$ws.Cells.Item(7, 4).Value = "CSC442"  # D7: None -> 'CSC442'
$ws.Cells.Item(7, 5).Value = "CSC442"  # E7: None -> 'CSC442'
$ws.Cells.Item(10, 6).Value = "BIO111"  # F10: None -> 'BIO111'
$ws.Cells.Item(13, 7).Value = ""  # G13: 'CIS421' -> None
$ws.Cells.Item(13, 11).Value = ""  # K13: 'CSC111' -> None
$ws.Cells.Item(16, 2).Value = "CSC424"  # B16: None -> 'CSC424'
$ws.Cells.Item(16, 3).Value = "CSC424"  # C16: None -> 'CSC424'
$ws.Cells.Item(19, 6).Value = "CIS421"  # F19: None -> 'CIS421'
$ws.Cells.Item(19, 7).Value = "CIS421"  # G19: None -> 'CIS421'
$ws.Cells.Item(20, 5).Value = "MAT112"  # E20: None -> 'MAT112'
$ws.Cells.Item(20, 6).Value = "MAT112"  # F20: None -> 'MAT112'
$ws.Cells.Item(20, 9).Value = "CST111"  # I20: 'TMC111' -> 'CST111'
$ws.Cells.Item(20, 10).Value = "CST111"  # J20: None -> 'CST111'
$ws.Cells.Item(21, 7).Value = "MAT111"  # G21: None -> 'MAT111'
$ws.Cells.Item(21, 8).Value = "MAT111"  # H21: None -> 'MAT111'
$ws.Cells.Item(21, 10).Value = ""  # J21: 'MAT112' -> None
$ws.Cells.Item(21, 11).Value = ""  # K21: 'MAT112' -> None
$ws.Cells.Item(24, 6).Value = ""  # F24: 'CSC425' -> None
$ws.Cells.Item(25, 4).Value = ""  # D25: 'CSC424' -> None
$ws.Cells.Item(25, 9).Value = "CSC424"  # I25: None -> 'CSC424'
$ws.Cells.Item(25, 10).Value = "CSC424"  # J25: None -> 'CSC424'

$ws = $wb.Worksheets.Item("thur")
$ws.Cells.Item(7, 5).Value = ""  # E7: 'CSC425' -> None
$ws.Cells.Item(7, 9).Value = "CSC111"  # I7: None -> 'CSC111'
$ws.Cells.Item(7, 10).Value = "CSC111"  # J7: None -> 'CSC111'
$ws.Cells.Item(11, 7).Value = ""  # G11: 'CSC424' -> None
$ws.Cells.Item(14, 9).Value = ""  # I14: 'CIS421' -> None
$ws.Cells.Item(16, 6).Value = "CSC425"  # F16: None -> 'CSC425'
$ws.Cells.Item(16, 7).Value = "CSC425"  # G16: None -> 'CSC425'
$ws.Cells.Item(17, 7).Value = "CSC442"  # G17: None -> 'CSC442'
$ws.Cells.Item(21, 2).Value = ""  # B21: 'EDS421' -> None
$ws.Cells.Item(21, 6).Value = "CHM111"  # F21: None -> 'CHM111'
$ws.Cells.Item(24, 11).Value = ""  # K24: 'CSC442' -> None
$ws.Cells.Item(25, 10).Value = ""  # J25: 'CSC424' -> None

$ws = $wb.Worksheets.Item("fri")
$ws.Cells.Item(2, 2).Value = ""  # B2: 'BIO111' -> None
$ws.Cells.Item(15, 5).Value = ""  # E15: 'CSC111' -> None
$ws.Cells.Item(17, 3).Value = ""  # C17: 'BIO111' -> None
$ws.Cells.Item(20, 3).Value = "MAT111"  # C20: None -> 'MAT111'
$ws.Cells.Item(20, 6).Value = "TMC111"  # F20: 'TMC421' -> 'TMC111'
$ws.Cells.Item(21, 5).Value = ""  # E21: 'GST111' -> None
$ws.Cells.Item(24, 3).Value = ""  # C24: 'CSC424' -> None
